$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '95.220.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.70%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.383.28'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.94%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.83'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '641.73'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.39%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.404'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.39%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.990'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.61%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.366.64'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.83'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +8.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.197'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '95.095.24'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.99'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.030.55'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000248'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.48'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.391.31'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.80'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +6.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.59'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +12.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.490'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +10.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '505.05'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.22'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000191'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.43'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '91.08'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.04%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.18'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.577.04'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.79'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +10.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.997'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.74'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +12.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.137'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.181'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.45%  '

$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.569'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.84%  '

$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '30.21'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +8.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.66'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.47%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.13%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.149'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '500.71'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.890'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +8.92%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.15'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.42%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0414'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.99%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.54%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.45'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.65%  '

$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.16'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +10.14%  '

$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.24'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.95%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.22'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.77%  '
